$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.115.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.836.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6817"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2988"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07658"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.021"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6776"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.163"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.22%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008287"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.037.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.068.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9983"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.341"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9991"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1441"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.716"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.502"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.253"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.132"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.197"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05403"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7511"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.859"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.130"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.300.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01815"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.715"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9360"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.052"
$ws.Range("D42").Style = "Normal"
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("B44").Value = "XinFinNetwork"
$ws.Range("C44").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.08273"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +29.92%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.966.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.765"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.371"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.17%  "
